# Update generated output values (想去人数 / "want-to-go" counts)
# for the 南宁-漫展信息 workbook, as produced by the gh-pages build at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): row for 熊喵M动漫嘉年华 (F3) and 第二届北极光动漫展 (F4)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1278
$wsExpo.Range("F4").Value = 2779

# Sheet "全部类型" (All types): same two events appear at F5 and F6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1278
$wsAll.Range("F6").Value = 2779
